$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 7606.143
$ws.Range("J18").Value = 100002
$ws.Range("L18").Value = 100002
$ws.Range("N18").Value = -100570
$ws.Range("H129").Value = 2873.6274
$ws.Range("J129").Value = 977.5217
$ws.Range("L129").Value = 2932.5651
$ws.Range("N129").Value = -12932.5651
$ws.Range("H130").Value = 43770
$ws.Range("J130").Value = 43770
$ws.Range("L130").Value = 43770
$ws.Range("N130").Value = -53810
$ws.Range("H132").Value = 5687393
$ws.Range("I132").Value = 5957718.5
$ws.Range("J132").Value = 10555.5
$ws.Range("K132").Value = 17873155.5
$ws.Range("L132").Value = 31666.5
$ws.Range("M132").Value = -17870625.5
$ws.Range("N132").Value = -36726.5
$ws.Range("H137").Value = 1551.1892
$ws.Range("I137").Value = 1189.6
$ws.Range("J137").Value = 3100.8572
$ws.Range("K137").Value = 3568.8
$ws.Range("L137").Value = 9302.571599999999
$ws.Range("M137").Value = -1018.8
$ws.Range("N137").Value = -14402.5716
$ws.Range("H141").Value = 1660.15
$ws.Range("I141").Value = 1555.7778
$ws.Range("J141").Value = 2599.5
$ws.Range("K141").Value = 4667.3334
$ws.Range("L141").Value = 7798.5
$ws.Range("M141").Value = 512.6665999999996
$ws.Range("N141").Value = -18158.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1441.2858
$ws.Range("I61").Value = 1346.6857
$ws.Range("J61").Value = 1914.2858
$ws.Range("K61").Value = 1346.6857
$ws.Range("L61").Value = 1914.2858
$ws.Range("M61").Value = -1134.6857
$ws.Range("N61").Value = -2338.2858
$ws.Range("H74").Value = 1847.091
$ws.Range("I74").Value = 827.91174
$ws.Range("J74").Value = 5312.3
$ws.Range("K74").Value = 827.91174
$ws.Range("L74").Value = 5312.3
$ws.Range("M74").Value = 46.08825999999999
$ws.Range("N74").Value = -7060.3
$ws.Range("H77").Value = 1847.091
$ws.Range("I77").Value = 827.91174
$ws.Range("J77").Value = 5312.3
$ws.Range("K77").Value = 4139.5587
$ws.Range("L77").Value = 26561.5
$ws.Range("M77").Value = 228.4413000000004
$ws.Range("N77").Value = -35297.5
$ws.Range("H132").Value = 2031.5778
$ws.Range("I132").Value = 1634.3448
$ws.Range("J132").Value = 2751.5625
$ws.Range("K132").Value = 4903.0344
$ws.Range("L132").Value = 8254.6875
$ws.Range("M132").Value = -2373.0344
$ws.Range("N132").Value = -13314.6875
$ws.Range("H136").Value = 1441.2858
$ws.Range("I136").Value = 1346.6857
$ws.Range("J136").Value = 1914.2858
$ws.Range("K136").Value = 4040.0571
$ws.Range("L136").Value = 5742.857400000001
$ws.Range("M136").Value = -1490.0571
$ws.Range("N136").Value = -10842.8574

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 86183.336
$ws.Range("I20").Value = 93827.27
$ws.Range("J20").Value = 2100
$ws.Range("K20").Value = 93827.27
$ws.Range("L20").Value = 2100
$ws.Range("M20").Value = -93580.27
$ws.Range("N20").Value = -2594
$ws.Range("H86").Value = 47012.668
$ws.Range("I86").Value = 59483.24
$ws.Range("J86").Value = 3365.6667
$ws.Range("K86").Value = 59483.24
$ws.Range("L86").Value = 3365.6667
$ws.Range("M86").Value = -58360.24
$ws.Range("N86").Value = -5611.6667
$ws.Range("H89").Value = 47012.668
$ws.Range("I89").Value = 59483.24
$ws.Range("J89").Value = 3365.6667
$ws.Range("K89").Value = 297416.2
$ws.Range("L89").Value = 16828.3335
$ws.Range("M89").Value = -291800.2
$ws.Range("N89").Value = -28060.3335
$ws.Range("H134").Value = 2686.7932
$ws.Range("I134").Value = 2864.12
$ws.Range("J134").Value = 1578.5
$ws.Range("K134").Value = 8592.360000000001
$ws.Range("L134").Value = 4735.5
$ws.Range("M134").Value = -6057.360000000001
$ws.Range("N134").Value = -9805.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 909.2857
$ws.Range("I22").Value = 405.2
$ws.Range("J22").Value = 1189.3334
$ws.Range("K22").Value = 405.2
$ws.Range("L22").Value = 1189.3334
$ws.Range("M22").Value = -55.19999999999999
$ws.Range("N22").Value = -1889.3334
$ws.Range("H31").Value = 43920.984
$ws.Range("I31").Value = 39769.58
$ws.Range("J31").Value = 47004.887
$ws.Range("K31").Value = 39769.58
$ws.Range("L31").Value = 47004.887
$ws.Range("M31").Value = -39474.58
$ws.Range("N31").Value = -47594.887
$ws.Range("H34").Value = 43920.984
$ws.Range("I34").Value = 39769.58
$ws.Range("J34").Value = 47004.887
$ws.Range("K34").Value = 39769.58
$ws.Range("L34").Value = 47004.887
$ws.Range("M34").Value = -39567.58
$ws.Range("N34").Value = -47408.887
$ws.Range("H58").Value = 6434.636
$ws.Range("I58").Value = 1455.8
$ws.Range("J58").Value = 14094.385
$ws.Range("K58").Value = 1455.8
$ws.Range("L58").Value = 14094.385
$ws.Range("M58").Value = -1252.8
$ws.Range("N58").Value = -14500.385
$ws.Range("H62").Value = 4275880.5
$ws.Range("I62").Value = 11113091
$ws.Range("J62").Value = 2623.5
$ws.Range("K62").Value = 11113091
$ws.Range("L62").Value = 2623.5
$ws.Range("M62").Value = -11112467
$ws.Range("N62").Value = -3871.5
$ws.Range("H65").Value = 4275880.5
$ws.Range("I65").Value = 11113091
$ws.Range("J65").Value = 2623.5
$ws.Range("K65").Value = 55565455
$ws.Range("L65").Value = 13117.5
$ws.Range("M65").Value = -55562335
$ws.Range("N65").Value = -19357.5
$ws.Range("H132").Value = 3713.6667
$ws.Range("I132").Value = 3369.111
$ws.Range("J132").Value = 4747.3335
$ws.Range("K132").Value = 10107.333
$ws.Range("L132").Value = 14242.0005
$ws.Range("M132").Value = -7577.332999999999
$ws.Range("N132").Value = -19302.0005
$ws.Range("H134").Value = 1248
$ws.Range("I134").Value = 1136.7693
$ws.Range("K134").Value = 3410.3079
$ws.Range("M134").Value = -875.3078999999998
$ws.Range("H136").Value = 6434.636
$ws.Range("I136").Value = 1455.8
$ws.Range("J136").Value = 14094.385
$ws.Range("K136").Value = 4367.4
$ws.Range("L136").Value = 42283.155
$ws.Range("M136").Value = -1817.4
$ws.Range("N136").Value = -47383.155

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 497
$ws.Range("J34").Value = 721.2
$ws.Range("L34").Value = 2163.6
$ws.Range("N34").Value = -2331.6
$ws.Range("H107").Value = 497041.25
$ws.Range("I107").Value = 527.0769
$ws.Range("J107").Value = 927353.5600000001
$ws.Range("K107").Value = 1581.2307
$ws.Range("L107").Value = 2782060.68
$ws.Range("M107").Value = 338.7692999999999
$ws.Range("N107").Value = -2785900.68
$ws.Range("H131").Value = 799.71
$ws.Range("J131").Value = 842.5402
$ws.Range("L131").Value = 2527.6206
$ws.Range("N131").Value = -12607.6206

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1504264.8
$ws.Range("I43").Value = 5000530.5
$ws.Range("J43").Value = 5865.2856
$ws.Range("K43").Value = 5000530.5
$ws.Range("L43").Value = 5865.2856
$ws.Range("M43").Value = -5000379.5
$ws.Range("N43").Value = -6167.2856
$ws.Range("H132").Value = 3107.625
$ws.Range("I132").Value = 1682.4
$ws.Range("K132").Value = 5047.200000000001
$ws.Range("M132").Value = -2517.200000000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1125628
$ws.Range("I46").Value = 490
$ws.Range("J46").Value = 1266270.2
$ws.Range("K46").Value = 490
$ws.Range("L46").Value = 1266270.2
$ws.Range("M46").Value = -302
$ws.Range("N46").Value = -1266646.2
$ws.Range("H68").Value = 2646
$ws.Range("I68").Value = 1668.4166
$ws.Range("J68").Value = 4321.857
$ws.Range("K68").Value = 1668.4166
$ws.Range("L68").Value = 4321.857
$ws.Range("M68").Value = -919.4166
$ws.Range("N68").Value = -5819.857
$ws.Range("H71").Value = 2646
$ws.Range("I71").Value = 1668.4166
$ws.Range("J71").Value = 4321.857
$ws.Range("K71").Value = 8342.083000000001
$ws.Range("L71").Value = 21609.285
$ws.Range("M71").Value = -4598.083000000001
$ws.Range("N71").Value = -29097.285
$ws.Range("H128").Value = 47932.5
$ws.Range("J128").Value = 47932.5
$ws.Range("L128").Value = 47932.5
$ws.Range("N128").Value = -57892.5
$ws.Range("H132").Value = 3238.4138
$ws.Range("I132").Value = 3318.087
$ws.Range("J132").Value = 2933
$ws.Range("K132").Value = 9954.261
$ws.Range("L132").Value = 8799
$ws.Range("M132").Value = -7424.261
$ws.Range("N132").Value = -13859
$ws.Range("H136").Value = 1035.579
$ws.Range("I136").Value = 996
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 2988
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -438
$ws.Range("N136").Value = -12600

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4275866
$ws.Range("J62").Value = 2609.0908
$ws.Range("L62").Value = 2609.0908
$ws.Range("N62").Value = -3857.0908
$ws.Range("H65").Value = 4275866
$ws.Range("J65").Value = 2609.0908
$ws.Range("L65").Value = 13045.454
$ws.Range("N65").Value = -19285.454
$ws.Range("H122").Value = 2007.5652
$ws.Range("I122").Value = 1045.4286
$ws.Range("J122").Value = 3504.2222
$ws.Range("K122").Value = 3136.2858
$ws.Range("L122").Value = 10512.6666
$ws.Range("M122").Value = -686.2857999999997
$ws.Range("N122").Value = -15412.6666
$ws.Range("H132").Value = 3121.5483
$ws.Range("I132").Value = 2955.1853
$ws.Range("J132").Value = 4244.5
$ws.Range("K132").Value = 8865.555899999999
$ws.Range("L132").Value = 12733.5
$ws.Range("M132").Value = -6335.555899999999
$ws.Range("N132").Value = -17793.5
$ws.Range("H136").Value = 1151.2307
$ws.Range("I136").Value = 870.75
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 2612.25
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -62.25
$ws.Range("N136").Value = -9900
